$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Boiler-1: drop the first two data rows (old rows 2 & 3). The old row 4
# (200.5 / -50.5 / 400 / 2000 / 25 / -50.5 / -50.5 / "-10%") shifts up to
# become the sole remaining data row (new row 2).
# ---------------------------------------------------------------------------
$wsBoiler1 = $wb.Worksheets.Item("Boiler-1")
$wsBoiler1.Range("A2:A3").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Boiler-2: insert a new row above the existing "N/A" row (pushing it from
# row 2 down to row 3), then fill the new row 2 with fresh numeric data.
# ---------------------------------------------------------------------------
$wsBoiler2 = $wb.Worksheets.Item("Boiler-2")
$wsBoiler2.Range("A2").EntireRow.Insert()

$wsBoiler2.Cells.Item(2, 1).Value = 187.05
$wsBoiler2.Cells.Item(2, 2).Value = 147.76
$wsBoiler2.Cells.Item(2, 3).Value = 307.69
$wsBoiler2.Cells.Item(2, 4).Value = 3835.14
$wsBoiler2.Cells.Item(2, 5).Value = 42.2
$wsBoiler2.Cells.Item(2, 6).Value = 147.76
$wsBoiler2.Cells.Item(2, 7).Value = 169.42
$wsBoiler2.Cells.Item(2, 8).Value = "'94.6%"

# ---------------------------------------------------------------------------
# Turbine-A: insert two new rows below the existing data row (old row 2
# stays put but gets overwritten; old row 3 — the "1,234.50" text row —
# shifts down from row 3 to row 5). Then populate rows 2-4 with fresh data.
# ---------------------------------------------------------------------------
$wsTurbineA = $wb.Worksheets.Item("Turbine-A")
$wsTurbineA.Range("A3:A4").EntireRow.Insert()

$wsTurbineA.Cells.Item(2, 1).Value = 293.67
$wsTurbineA.Cells.Item(2, 2).Value = 156.4
$wsTurbineA.Cells.Item(2, 3).Value = 303.3
$wsTurbineA.Cells.Item(2, 4).Value = 2403.29
$wsTurbineA.Cells.Item(2, 5).Value = 25.18
$wsTurbineA.Cells.Item(2, 6).Value = 156.4
$wsTurbineA.Cells.Item(2, 7).Value = 123.09
$wsTurbineA.Cells.Item(2, 8).Value = "'86.2%"

$wsTurbineA.Cells.Item(3, 1).Value = 459.68
$wsTurbineA.Cells.Item(3, 2).Value = 50.75
$wsTurbineA.Cells.Item(3, 3).Value = 426.35
$wsTurbineA.Cells.Item(3, 4).Value = 1795.1
$wsTurbineA.Cells.Item(3, 5).Value = 29.14
$wsTurbineA.Cells.Item(3, 6).Value = 50.75
$wsTurbineA.Cells.Item(3, 7).Value = 190.13
$wsTurbineA.Cells.Item(3, 8).Value = "'88.7%"

$wsTurbineA.Cells.Item(4, 1).Value = 447.51
$wsTurbineA.Cells.Item(4, 2).Value = 84.31999999999999
$wsTurbineA.Cells.Item(4, 3).Value = 449.87
$wsTurbineA.Cells.Item(4, 4).Value = 2114.78
$wsTurbineA.Cells.Item(4, 5).Value = 36.42
$wsTurbineA.Cells.Item(4, 6).Value = 84.31999999999999
$wsTurbineA.Cells.Item(4, 7).Value = 143.63
$wsTurbineA.Cells.Item(4, 8).Value = "'91.7%"

# Row 5 (shifted down from the old row 3) already carries the correct
# "1,234.50" / "150 MT" / "450,0" / "2,500.99" / "35.50" / "150 MT" /
# "150 MT" / "88.5%" text values untouched, so nothing else to do here.

# ---------------------------------------------------------------------------
# Cooling-Tower: remove the old row 3 (225.48 data row) so the old row 4
# (105 / 60 / 320 / 4500 / 15 / 60 / 65 / "92%") shifts up to row 3, then
# overwrite row 2 with fresh data.
# ---------------------------------------------------------------------------
$wsCoolingTower = $wb.Worksheets.Item("Cooling-Tower")
$wsCoolingTower.Range("A3").EntireRow.Delete()

$wsCoolingTower.Cells.Item(2, 1).Value = 336.46
$wsCoolingTower.Cells.Item(2, 2).Value = 149.1
$wsCoolingTower.Cells.Item(2, 3).Value = 347.61
$wsCoolingTower.Cells.Item(2, 4).Value = 1985.8
$wsCoolingTower.Cells.Item(2, 5).Value = 26.09
$wsCoolingTower.Cells.Item(2, 6).Value = 149.1
$wsCoolingTower.Cells.Item(2, 7).Value = 228.91
$wsCoolingTower.Cells.Item(2, 8).Value = "'89.5%"
